# Auto-generated Word COM-interop script applying the target diff.
$d = $word.ActiveDocument

# --- Part 1: consolidate split runs / drop grammar-check proofErr markers. ---
# Each paragraph's text is unchanged; Find/Replace with the identical text
# re-merges the run(s) it spans into a single run and drops the <w:proofErr/>
# wrappers that previously surrounded the grammar-checked word.

$null = $d.Content.Find.Execute('Be lehet állítani, hogy a Drawer az AppBart eltolja vagy fölé gördüljön, illetve ', $true, $false, $false, $false, $false, $true, 1, $false, 'Be lehet állítani, hogy a Drawer az AppBart eltolja vagy fölé gördüljön, illetve ', 2)
$null = $d.Content.Find.Execute('yomására jelenjenek meg (pl.: egy Person típust megjelenítő sorhoz, az Address típust megjelenítő táblázat mint ChildRowContent)', $true, $false, $false, $false, $false, $true, 1, $false, 'yomására jelenjenek meg (pl.: egy Person típust megjelenítő sorhoz, az Address típust megjelenítő táblázat mint ChildRowContent)', 2)
$null = $d.Content.Find.Execute('Lézetik külön <MudAutoComplete> komponens is, ami amellett, hogy tud hasonlóan viselkedni mint a <MudSelect> komponens, képes egyéni bemenetet is kezelni, illetve meg lehet neki adni egy search function-t, aminek segítségével autocomplete viselkedést valósít meg.', $true, $false, $false, $false, $false, $true, 1, $false, 'Lézetik külön <MudAutoComplete> komponens is, ami amellett, hogy tud hasonlóan viselkedni mint a <MudSelect> komponens, képes egyéni bemenetet is kezelni, illetve meg lehet neki adni egy search function-t, aminek segítségével autocomplete viselkedést valósít meg.', 2)
$null = $d.Content.Find.Execute('Van, a <MudProgressCircular> és <MudProgressLinear> komponensekkel. A komponens lehet Determinate, ha meg lehet becsülni hol tart egy adott folyamat, illetve Indeterminate ha csak azt akarjuk mutatni, hogy valamilyen háttérfolyamatra, de nem tudjuk az hol tart.', $true, $false, $false, $false, $false, $true, 1, $false, 'Van, a <MudProgressCircular> és <MudProgressLinear> komponensekkel. A komponens lehet Determinate, ha meg lehet becsülni hol tart egy adott folyamat, illetve Indeterminate ha csak azt akarjuk mutatni, hogy valamilyen háttérfolyamatra, de nem tudjuk az hol tart.', 2)
$null = $d.Content.Find.Execute(' A <MudForm>-ot legegyszerűbben úgy lehet validálni, hogy az egyes bemeneti mezőkhöz megadjuk a Required illetve Validation attribútumokat. A <MudForm> rendelkezik egy bind-IsValid attribútummal, amivel a validáció sikerességét tudjuk változóba menteni, illetve egy bind-Errors attribútuma, amivel a hiba stringeket tudjuk összegyűjteni. ', $true, $false, $false, $false, $false, $true, 1, $false, ' A <MudForm>-ot legegyszerűbben úgy lehet validálni, hogy az egyes bemeneti mezőkhöz megadjuk a Required illetve Validation attribútumokat. A <MudForm> rendelkezik egy bind-IsValid attribútummal, amivel a validáció sikerességét tudjuk változóba menteni, illetve egy bind-Errors attribútuma, amivel a hiba stringeket tudjuk összegyűjteni. ', 2)
$null = $d.Content.Find.Execute('Milyen grafikonok vannak? (kördiagram, oszlopdiagram, stb.)', $true, $false, $false, $false, $false, $true, 1, $false, 'Milyen grafikonok vannak? (kördiagram, oszlopdiagram, stb.)', 2)
$null = $d.Content.Find.Execute('<MudAlert>: Egyszerű kis komponens, ami egy ikont és kis szöveget tud megjeleníteni. A bordert és a színeket testre lehet szabni. Alapvetően statikus, de a sarkában lehet megjeleníteni kis X-et amivel bezárhatjuk, azonban ezt pár sor kóddal kell megoldanunk.', $true, $false, $false, $false, $false, $true, 1, $false, '<MudAlert>: Egyszerű kis komponens, ami egy ikont és kis szöveget tud megjeleníteni. A bordert és a színeket testre lehet szabni. Alapvetően statikus, de a sarkában lehet megjeleníteni kis X-et amivel bezárhatjuk, azonban ezt pár sor kóddal kell megoldanunk.', 2)
$null = $d.Content.Find.Execute('<MudSnackbarProvider>: Hasonlóan a MudDialogProviderhez, ezt is érdemes globálisan egyszer felvenni, majd injektálni ahol kell. A Snackbart, az Alert-tel ellentétben, kódból dinamikusan tudjuk csak megjeleníteni a SnackbarProvider Add metódusával. A tartalma lehet ikon, szöveg, action button, de akár még custom HTML markup is.', $true, $false, $false, $false, $false, $true, 1, $false, '<MudSnackbarProvider>: Hasonlóan a MudDialogProviderhez, ezt is érdemes globálisan egyszer felvenni, majd injektálni ahol kell. A Snackbart, az Alert-tel ellentétben, kódból dinamikusan tudjuk csak megjeleníteni a SnackbarProvider Add metódusával. A tartalma lehet ikon, szöveg, action button, de akár még custom HTML markup is.', 2)

# --- Part 2: new file-upload documentation content. ---

# Locate the two placeholder (empty) list paragraphs that sit right below
# the "Van-e lehetőség egy vagy több fájl feltöltésére?" / "...hol tart a
# feltöltés?" headings, by scanning for paragraphs whose text is empty.
$targets = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Length -le 1) {
        $targets += $p
    }
}

# The last two empty list paragraphs in the document are the two we need
# (the upload-progress Q&A block at the very end of the document).
$uploadPara = $targets[$targets.Count - 2]
$progressPara = $targets[$targets.Count - 1]

# Fill the first placeholder with the <InputFile> explanation AND insert a
# brand-new paragraph right after it (OnChanged / InputFileChangeEventArgs),
# in one InsertXML call so both paragraphs land with no stray paraId/rsid.
$null = $uploadPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listaszerbekezds"/></w:pPr><w:r><w:t xml:space="preserve">Van, a natív &lt;</w:t></w:r><w:r><w:t xml:space="preserve">InputFile</w:t></w:r><w:r><w:t xml:space="preserve">&gt; komponenst felhasználva. Ha bármilyen MudButton for attribútumát beállítjuk az input id attribútumának értékére. Ennek a gombnak aztán olyan stílust adhatunk amilyet csak szeretnénk, de akár másféle komponenst is használhatunk, ami képes Click event kezelésére. Akár Drag-and-Drop fájlfeltöltést is meg tudunk valósítani, ha olyan komponenst használunk az input label-jének, ami tudja kezelni a drag eventeket.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listaszerbekezds"/></w:pPr><w:r><w:t xml:space="preserve">Magát a feltöltés kezdeményezését az input OnChanged eventjével tudjuk kezelni. Az </w:t></w:r><w:r><w:t xml:space="preserve">InputFileChangeEventArgs</w:t></w:r><w:r><w:t xml:space="preserve"> paraméteren keresztül tudjuk bekérni a fájlokat.</w:t></w:r></w:p>')

# Fill the second placeholder with the upload-progress answer.
$null = $progressPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listaszerbekezds"/></w:pPr><w:r><w:t xml:space="preserve">Nincs a könyvtárnak beépített megoldása erre. Saját kézzel úgy oldhatjuk meg ezt, hogy chunkonként olvassuk be a fájlt, és valamilyen Timer segítségével frissítjük a progress mérőt.</w:t></w:r></w:p>')

